$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 880.9
$ws.Range("I15").Value = 880.9
$ws.Range("K15").Value = 2642.7
$ws.Range("M15").Value = -2473.7
$ws.Range("H33").Value = 299
$ws.Range("I33").Value = 299
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 299
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -70
$ws.Range("N33").ClearContents()
$ws.Range("H34").Value = 3371.4285
$ws.Range("I34").Value = 3371.4285
$ws.Range("K34").Value = 3371.4285
$ws.Range("M34").Value = -3168.4285
$ws.Range("H36").Value = 3371.4285
$ws.Range("I36").Value = 3371.4285
$ws.Range("K36").Value = 3371.4285
$ws.Range("M36").Value = -2656.4285
$ws.Range("H62").Value = 4048.5
$ws.Range("I62").Value = 3069.2856
$ws.Range("J62").Value = 6333.3335
$ws.Range("K62").Value = 3069.2856
$ws.Range("L62").Value = 6333.3335
$ws.Range("M62").Value = -2445.2856
$ws.Range("N62").Value = -7581.3335
$ws.Range("H65").Value = 4048.5
$ws.Range("I65").Value = 3069.2856
$ws.Range("J65").Value = 6333.3335
$ws.Range("K65").Value = 15346.428
$ws.Range("L65").Value = 31666.6675
$ws.Range("M65").Value = -12226.428
$ws.Range("N65").Value = -37906.6675
$ws.Range("H95").Value = 40000
$ws.Range("J95").Value = 40000
$ws.Range("L95").Value = 40000
$ws.Range("N95").Value = -45492
$ws.Range("H98").Value = 4280.8696
$ws.Range("I98").Value = 1678.1818
$ws.Range("K98").Value = 1678.1818
$ws.Range("M98").Value = -180.1818000000001
$ws.Range("H122").Value = 4280.8696
$ws.Range("I122").Value = 1678.1818
$ws.Range("K122").Value = 5034.5454
$ws.Range("M122").Value = -2584.5454
$ws.Range("H128").Value = 41846
$ws.Range("J128").Value = 41846
$ws.Range("L128").Value = 41846
$ws.Range("N128").Value = -51806
$ws.Range("H129").Value = 855.08
$ws.Range("I129").Value = 366.66666
$ws.Range("J129").Value = 870.18555
$ws.Range("K129").Value = 1099.99998
$ws.Range("L129").Value = 2610.55665
$ws.Range("M129").Value = 3900.00002
$ws.Range("N129").Value = -12610.55665
$ws.Range("H132").Value = 23257380
$ws.Range("I132").Value = 24391398
$ws.Range("K132").Value = 73174194
$ws.Range("M132").Value = -73171664
$ws.Range("H138").Value = 2704.293
$ws.Range("I138").Value = 941.7857
$ws.Range("J138").Value = 2994.5881
$ws.Range("K138").Value = 2825.3571
$ws.Range("L138").Value = 8983.764299999999
$ws.Range("M138").Value = 2314.6429
$ws.Range("N138").Value = -19263.7643
$ws.Range("H141").Value = 135252.47
$ws.Range("I141").Value = 155306.69
$ws.Range("K141").Value = 465920.07
$ws.Range("M141").Value = -460740.07

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3930.077
$ws.Range("I32").Value = 3122.2131
$ws.Range("J32").Value = 16250
$ws.Range("K32").Value = 3122.2131
$ws.Range("L32").Value = 16250
$ws.Range("M32").Value = -2835.2131
$ws.Range("N32").Value = -16824
$ws.Range("H45").Value = 1452
$ws.Range("I45").Value = 1526.6666
$ws.Range("J45").Value = 1340
$ws.Range("K45").Value = 1526.6666
$ws.Range("L45").Value = 1340
$ws.Range("M45").Value = -1149.6666
$ws.Range("N45").Value = -2094
$ws.Range("H47").Value = 36498.5
$ws.Range("J47").Value = 36498.5
$ws.Range("L47").Value = 36498.5
$ws.Range("N47").Value = -37948.5
$ws.Range("H122").Value = 7026.2856
$ws.Range("I122").Value = 1658.5
$ws.Range("J122").Value = 14183.333
$ws.Range("K122").Value = 4975.5
$ws.Range("L122").Value = 42549.999
$ws.Range("M122").Value = -2525.5
$ws.Range("N122").Value = -47449.999
$ws.Range("H132").Value = 1921.3438
$ws.Range("I132").Value = 996.9545000000001
$ws.Range("J132").Value = 3955
$ws.Range("K132").Value = 2990.8635
$ws.Range("L132").Value = 11865
$ws.Range("M132").Value = -460.8635000000004
$ws.Range("N132").Value = -16925

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H34").Value = 2222
$ws.Range("I34").Value = 2222
$ws.Range("K34").Value = 2222
$ws.Range("M34").Value = -2108
$ws.Range("H42").Value = 79800
$ws.Range("J42").Value = 79800
$ws.Range("L42").Value = 79800
$ws.Range("N42").Value = -80456
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H105").Value = 1649.2222
$ws.Range("I105").Value = 1626.8358
$ws.Range("J105").Value = 1949.2
$ws.Range("K105").Value = 1626.8358
$ws.Range("L105").Value = 1949.2
$ws.Range("M105").Value = 120.1641999999999
$ws.Range("N105").Value = -5443.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2963.742
$ws.Range("I31").Value = 964.26666
$ws.Range("J31").Value = 4838.25
$ws.Range("K31").Value = 964.26666
$ws.Range("L31").Value = 4838.25
$ws.Range("M31").Value = -669.26666
$ws.Range("N31").Value = -5428.25
$ws.Range("H34").Value = 2963.742
$ws.Range("I34").Value = 964.26666
$ws.Range("J34").Value = 4838.25
$ws.Range("K34").Value = 964.26666
$ws.Range("L34").Value = 4838.25
$ws.Range("M34").Value = -762.26666
$ws.Range("N34").Value = -5242.25
$ws.Range("H122").Value = 3553
$ws.Range("I122").Value = 1555.5
$ws.Range("J122").Value = 4218.8335
$ws.Range("K122").Value = 4666.5
$ws.Range("L122").Value = 12656.5005
$ws.Range("M122").Value = -2216.5
$ws.Range("N122").Value = -17556.5005
$ws.Range("H137").Value = 45438.57
$ws.Range("J137").Value = 45438.57
$ws.Range("L137").Value = 45438.57
$ws.Range("N137").Value = -55638.57

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 461614.44
$ws.Range("I5").Value = 623.7143
$ws.Range("J5").Value = 608293.3
$ws.Range("K5").Value = 1871.1429
$ws.Range("L5").Value = 1824879.9
$ws.Range("M5").Value = -1759.1429
$ws.Range("N5").Value = -1825103.9
$ws.Range("H113").Value = 3049373.2
$ws.Range("I113").Value = 600.1739
$ws.Range("K113").Value = 1800.5217
$ws.Range("M113").Value = 369.4783
$ws.Range("H122").Value = 2964.1316
$ws.Range("I122").Value = 1066.2727
$ws.Range("K122").Value = 9596.454299999999
$ws.Range("M122").Value = -7146.454299999999
$ws.Range("H132").Value = 1598.5186
$ws.Range("I132").Value = 775
$ws.Range("J132").Value = 3951.4285
$ws.Range("K132").Value = 6975
$ws.Range("L132").Value = 35562.8565
$ws.Range("M132").Value = -4445
$ws.Range("N132").Value = -40622.8565
$ws.Range("H135").Value = 461614.44
$ws.Range("I135").Value = 623.7143
$ws.Range("J135").Value = 608293.3
$ws.Range("K135").Value = 5613.428699999999
$ws.Range("L135").Value = 5474639.7
$ws.Range("M135").Value = -3078.428699999999
$ws.Range("N135").Value = -5479709.7
$ws.Range("H137").Value = 2772.8572
$ws.Range("I137").Value = 1427.5
$ws.Range("J137").Value = 4566.6665
$ws.Range("K137").Value = 4282.5
$ws.Range("L137").Value = 13699.9995
$ws.Range("M137").Value = 817.5
$ws.Range("N137").Value = -23899.9995
$ws.Range("H140").Value = 3852.2307
$ws.Range("I140").Value = 4986.5557
$ws.Range("J140").Value = 1300
$ws.Range("K140").Value = 14959.6671
$ws.Range("L140").Value = 3900
$ws.Range("M140").Value = -9779.667099999999
$ws.Range("N140").Value = -14260

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1927.4849
$ws.Range("I102").Value = 1247.6666
$ws.Range("K102").Value = 1247.6666
$ws.Range("M102").Value = 374.3334
$ws.Range("H122").Value = 8426
$ws.Range("I122").Value = 2900
$ws.Range("K122").Value = 8700
$ws.Range("M122").Value = -6250

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6878.737
$ws.Range("I40").Value = 5981.4546
$ws.Range("J40").Value = 8112.5
$ws.Range("K40").Value = 5981.4546
$ws.Range("L40").Value = 8112.5
$ws.Range("M40").Value = -5845.4546
$ws.Range("N40").Value = -8384.5
$ws.Range("H68").Value = 668.6869
$ws.Range("I68").Value = 668.6869
$ws.Range("K68").Value = 668.6869
$ws.Range("M68").Value = 80.31309999999996
$ws.Range("H71").Value = 668.6869
$ws.Range("I71").Value = 668.6869
$ws.Range("K71").Value = 3343.4345
$ws.Range("M71").Value = 400.5654999999997
$ws.Range("H86").Value = 34130
$ws.Range("J86").Value = 34130
$ws.Range("L86").Value = 34130
$ws.Range("N86").Value = -36502
$ws.Range("H89").Value = 34130
$ws.Range("J89").Value = 34130
$ws.Range("L89").Value = 102390
$ws.Range("N89").Value = -114246
$ws.Range("H122").Value = 8278.571
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 9241.666999999999
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 27725.001
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -32625.001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 39266.668
$ws.Range("J80").Value = 39266.668
$ws.Range("L80").Value = 39266.668
$ws.Range("N80").Value = -41262.668
$ws.Range("H83").Value = 39266.668
$ws.Range("J83").Value = 39266.668
$ws.Range("L83").Value = 117800.004
$ws.Range("N83").Value = -127784.004
$ws.Range("H122").Value = 3532.55
$ws.Range("I122").Value = 1493.6923
$ws.Range("J122").Value = 7319
$ws.Range("K122").Value = 4481.0769
$ws.Range("L122").Value = 21957
$ws.Range("M122").Value = -2031.0769
$ws.Range("N122").Value = -26857
$ws.Range("H132").Value = 13335881
$ws.Range("I132").Value = 1842.2858
$ws.Range("J132").Value = 30306476
$ws.Range("K132").Value = 5526.857400000001
$ws.Range("L132").Value = 90919428
$ws.Range("M132").Value = -2996.857400000001
$ws.Range("N132").Value = -90924488
